$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Evaporator Temperature
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 20.12729329293731
$ws.Range("D2").Value = 1.067392729962487
$ws.Range("F2").Value = 19.40741293415048
$ws.Range("G2").Value = 20.12745645159733
$ws.Range("H2").Value = 20.84707495745016

# Row 3: Condenser Temperature
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 49.99984
$ws.Range("D3").Value = 5.8028122923717

# Row 4: Adiabatic Efficiency
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 74.99988
$ws.Range("D4").Value = 2.91582242080651

# Row 5: Compressor Energy
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = 747.0884972501755
$ws.Range("D5").Value = 256.5742921706667
$ws.Range("E5").Value = 264.1617093327257
$ws.Range("F5").Value = 553.0671495963257
$ws.Range("G5").Value = 706.107070080722
$ws.Range("H5").Value = 909.764906974582
$ws.Range("I5").Value = 1610.199657421221

# Row 6: Electric Current
$ws.Range("B6").Value = 5000
$ws.Range("C6").Value = 3.395856805682616
$ws.Range("D6").Value = 1.16624678259394
$ws.Range("E6").Value = 1.20073504242148
$ws.Range("F6").Value = 2.513941589074208
$ws.Range("G6").Value = 3.209577591276009
$ws.Range("H6").Value = 4.135295031702645
$ws.Range("I6").Value = 7.319089351914639

# Row 7: Discharge Temperature
$ws.Range("B7").Value = 5000
$ws.Range("C7").Value = 66.59293258773384
$ws.Range("D7").Value = 8.829039347566445
$ws.Range("E7").Value = 49.37996184290603
$ws.Range("F7").Value = 58.98508752900509
$ws.Range("G7").Value = 66.58805403244403
$ws.Range("H7").Value = 74.29700139325315
$ws.Range("I7").Value = 84.07697108049842

# Row 8: Refrigerant Mass Flow
$ws.Range("B8").Value = 5000
$ws.Range("C8").Value = 1.715485958031988
$ws.Range("D8").Value = 0.3737026169477889
$ws.Range("E8").Value = 0.9945388127932807
$ws.Range("F8").Value = 1.413521339771791
$ws.Range("G8").Value = 1.702474104568811
$ws.Range("H8").Value = 1.989265083471993
$ws.Range("I8").Value = 2.67370787602378

# Row 9: Capacity
$ws.Range("B9").Value = 5000
$ws.Range("D9").Value = 2655.699257768007
$ws.Range("F9").Value = 11200
$ws.Range("H9").Value = 15800
